# Weekly update: a new price record is inserted at the top of the
# "Feria Lagunitas de Puerto Montt" / Espárragos data block (row 35),
# pushing the existing rows 35-48 down to 36-49.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 35:48 down to 36:49, leaving a blank row 35.
$ws.Rows("35:35").Insert()

# Populate the new row 35 with this week's record.
$ws.Range("A35").Value = 4
$ws.Range("B35").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C35").Value = "Los Lagos"
$ws.Range("D35").Value = 44845
$ws.Range("E35").Value = 10
$ws.Range("F35").Value = 300000000
$ws.Range("G35").Value = "Espárragos"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 400
$ws.Range("K35").Value = 1800
$ws.Range("L35").Value = 1800
$ws.Range("M35").Value = 1800
$ws.Range("N35").Value = "$/kilo"
$ws.Range("O35").Value = "Provincia de Linares"
$ws.Range("P35").Value = 1800
$ws.Range("Q35").Value = 1
$ws.Range("R35").Value = "Hortaliza"
